# "Generate Report for Handback"
# A new handback report pass re-stamped the handoff/handback timestamps
# for the file "11620a51-59f4-49f4-9da4-f99daa341c5b" (the row whose
# status is "Handed back: in sync with en-US") across the Overview sheet
# and each per-locale sheet. The sibling file
# "34592a2d-0df0-442a-9a09-29b8a21fef50" row is untouched.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
#     first data row (11620a51...) moves forward.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-04 00:53:20"

# --- zh-cn sheet: "Correspond Handoff Datetime" (H) / "Correspond
#     Handback DateTime" (K) for the first data row (11620a51...).
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-04 00:53:15"
$zhcn.Range("K2").Value = "2016-09-04 00:53:33"

# --- de-de sheet: same two columns, same row.
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-04 00:53:20"
$dede.Range("K2").Value = "2016-09-04 00:53:40"
